$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two "2025-07-21" rows (CARREGADOR USB-C A GOLD 20W CA31-4 and
# MARMITA ELETRICA ONEX) - these were rows 2 and 3.
$ws.Range("A2:A3").EntireRow.Delete()

# Update estoque_atualizado (column G) values for the refreshed data.
$ws.Range("G2").Value = -252
$ws.Range("G5").Value = -74
$ws.Range("G7").Value = -61
$ws.Range("G9").Value = -1251
